$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.621.57'
Set-TextValue 'E2' '  -0.31%  '
Set-TextValue 'D3' '1.595.65'
Set-TextValue 'E3' '  -0.30%  '
Set-TextValue 'E4' '  -0.03%  '
Set-TextValue 'D5' '210.69'
Set-TextValue 'E5' '  -0.34%  '
Set-TextValue 'E6' '  -0.53%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'E8' '  -0.64%  '
Set-TextValue 'E9' '  -0.47%  '
Set-TextValue 'D10' '19.58'
Set-TextValue 'E10' '  +0.04%  '
Set-TextValue 'D11' '0.0845'
Set-TextValue 'E11' '  +0.13%  '
Set-TextValue 'D12' '1.819.61'
Set-TextValue 'E12' '  -0.30%  '
Set-TextValue 'D13' '1.593.91'
Set-TextValue 'E13' '  -0.72%  '
Set-TextValue 'E14' '  -0.21%  '
Set-TextValue 'D15' '0.522'
Set-TextValue 'E15' '  -0.18%  '
Set-TextValue 'D16' '64.54'
Set-TextValue 'E16' '  -1.20%  '
Set-TextValue 'D17' '26.606.65'
Set-TextValue 'E18' '  -2.39%  '
Set-TextValue 'E19' '  -0.03%  '
Set-TextValue 'D20' '208.73'
Set-TextValue 'E20' '  -0.39%  '
Set-TextValue 'D21' '7.06'
Set-TextValue 'E21' '  -2.19%  '
Set-TextValue 'D22' '4.29'
Set-TextValue 'E22' '  +0.18%  '
Set-TextValue 'E23' '  -3.67%  '
Set-TextValue 'D25' '145.04'
Set-TextValue 'E25' '  +1.43%  '
Set-TextValue 'E26' '  -0.02%  '
Set-TextValue 'E27' '  +0.09%  '
Set-TextValue 'E28' '  -0.93%  '
Set-TextValue 'D29' '15.26'
Set-TextValue 'E29' '  -0.57%  '
Set-TextValue 'E30' '  -2.82%  '
Set-TextValue 'E31' '  -0.64%  '
Set-TextValue 'E32' '  -0.20%  '
Set-TextValue 'E33' '  -0.38%  '
Set-TextValue 'D34' '1.281.47'
Set-TextValue 'E34' '  -1.01%  '
Set-TextValue 'E35' '  +0.38%  '
Set-TextValue 'D36' '1.21'
Set-TextValue 'E36' '  +11.31%  '
Set-TextValue 'B37' 'LidoDAOToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D37' '1.48'
Set-TextValue 'E37' '  -1.29%  '
Set-TextValue 'B38' 'ImmutableX'
Set-TextValue 'C38' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D38' '0.600'
Set-TextValue 'E38' '  -3.29%  '
Set-TextValue 'D39' '0.0168'
Set-TextValue 'E39' '  -2.00%  '
Set-TextValue 'D40' '0.823'
Set-TextValue 'E40' '  -0.46%  '
Set-TextValue 'E42' '  -1.27%  '
Set-TextValue 'D43' '0.772'
Set-TextValue 'E43' '  -1.61%  '
Set-TextValue 'D44' '62.79'
Set-TextValue 'E44' '  -0.97%  '
Set-TextValue 'D45' '1.731.58'
Set-TextValue 'E45' '  -0.29%  '
Set-TextValue 'D46' '89.49'
Set-TextValue 'E46' '  -1.70%  '
Set-TextValue 'D47' '1.58'
Set-TextValue 'E47' '  +0.12%  '
Set-TextValue 'E48' '  +2.24%  '
Set-TextValue 'D49' '0.0512'
Set-TextValue 'E49' '  +0.57%  '
Set-TextValue 'B50' 'EnergySwap'
Set-TextValue 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '7.47'
Set-TextValue 'E50' '  +0.85%  '
Set-TextValue 'B51' 'USDD'
Set-TextValue 'C51' 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D51' '1.00'
Set-TextValue 'E51' '  -0.15%  '
